$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Equities")
$ws2 = $wb.Worksheets.Item("Commodities")

# --- Equities sheet: refresh data rows 2-39 ---
$ws1.Range("B2").Value = "China"
$ws1.Range("C2").Value = 13245.09
$ws1.Range("D2").Value = 302.13999999999942
$ws1.Range("E2").Value = 0.023343982631471191
$ws1.Range("F2").Value = 0.23805905066894001
$ws1.Range("G2").Value = 0.033813274325968923
$ws1.Range("H2").Value = 0.27473899163978333
$ws1.Range("B3").Value = "China"
$ws1.Range("C3").Value = 15678.89
$ws1.Range("D3").Value = 324.30999999999949
$ws1.Range("E3").Value = 0.021121385280483022
$ws1.Range("F3").Value = 0.093292787522191611
$ws1.Range("G3").Value = 0.031567938755592273
$ws1.Range("H3").Value = 0.1256837424515993
$ws1.Range("B4").Value = "United States"
$ws1.Range("C4").Value = 1536.78
$ws1.Range("D4").Value = 39.509999999999991
$ws1.Range("E4").Value = 0.02638802620769809
$ws1.Range("F4").Value = -0.076043432776595421
$ws1.Range("G4").Value = 0.02638802620769809
$ws1.Range("H4").Value = -0.076043432776595421
$ws1.Range("B5").Value = "Saudi Arabia"
$ws1.Range("C5").Value = 8337.8799999999992
$ws1.Range("D5").Value = 202.71999999999929
$ws1.Range("E5").Value = 0.02491899360307603
$ws1.Range("F5").Value = 0.016055071026171989
$ws1.Range("G5").Value = 0.024992928015268309
$ws1.Range("H5").Value = 0.017316173592687441
$ws1.Range("B6").Value = "Vietnam"
$ws1.Range("C6").Value = 239.64
$ws1.Range("D6").Value = 5.3299999999999841
$ws1.Range("E6").Value = 0.022747642012718169
$ws1.Range("F6").Value = 0.34803397648647111
$ws1.Range("G6").Value = 0.022879990589823599
$ws1.Range("H6").Value = 0.34809215675971039
$ws1.Range("B7").Value = "South Korea"
$ws1.Range("C7").Value = 2412.4
$ws1.Range("D7").Value = 15.71000000000004
$ws1.Range("E7").Value = 0.0065548735965017526
$ws1.Range("F7").Value = 0.1194067942108608
$ws1.Range("G7").Value = 0.02255772375300635
$ws1.Range("H7").Value = 0.1167247152468489
$ws1.Range("B8").Value = "Taiwan"
$ws1.Range("C8").Value = 12875.62
$ws1.Range("D8").Value = 199.6700000000001
$ws1.Range("E8").Value = 0.015751876585186949
$ws1.Range("F8").Value = 0.077154875281929103
$ws1.Range("G8").Value = 0.020861730567010461
$ws1.Range("H8").Value = 0.11092076643918761
$ws1.Range("B9").Value = "Chile"
$ws1.Range("C9").Value = 3729.64
$ws1.Range("D9").Value = 16.9699999999998
$ws1.Range("E9").Value = 0.0045708344668391474
$ws1.Range("F9").Value = -0.23641642047420461
$ws1.Range("G9").Value = 0.020285762122180358
$ws1.Range("H9").Value = -0.23330466323238769
$ws1.Range("B10").Value = "Denmark"
$ws1.Range("C10").Value = 1491.36
$ws1.Range("D10").Value = 26.029999999999969
$ws1.Range("E10").Value = 0.01776391666041088
$ws1.Range("F10").Value = 0.20742252015933141
$ws1.Range("G10").Value = 0.01987903691554704
$ws1.Range("H10").Value = 0.28738544676885791
$ws1.Range("B11").Value = "South Africa"
$ws1.Range("C11").Value = 50399.16
$ws1.Range("D11").Value = -1315.9599999999989
$ws1.Range("E11").Value = -0.025446329816115609
$ws1.Range("F11").Value = -0.011773619445325069
$ws1.Range("G11").Value = 0.01734196391073373
$ws1.Range("H11").Value = -0.1227815109654001
$ws1.Range("B12").Value = "Mexico"
$ws1.Range("C12").Value = 36017.35
$ws1.Range("D12").Value = -317.54000000000087
$ws1.Range("E12").Value = -0.0087392586024066521
$ws1.Range("F12").Value = -0.1905358543486616
$ws1.Range("G12").Value = 0.017066237679618901
$ws1.Range("H12").Value = -0.26600509222679197
$ws1.Range("B13").Value = "Russia"
$ws1.Range("C13").Value = 2951.79
$ws1.Range("D13").Value = 41.279999999999752
$ws1.Range("E13").Value = 0.014183081315645611
$ws1.Range("F13").Value = -0.041274883317580779
$ws1.Range("G13").Value = 0.016703773116527602
$ws1.Range("H13").Value = -0.2066177333281812
$ws1.Range("B14").Value = "South Korea"
$ws1.Range("C14").Value = 888.88
$ws1.Range("D14").Value = 0.43999999999994088
$ws1.Range("E14").Value = 0.00049525010130113678
$ws1.Range("F14").Value = 0.35642672933420833
$ws1.Range("G14").Value = 0.016401760505903251
$ws1.Range("H14").Value = 0.3531767547800202
$ws1.Range("B15").Value = "Brazil"
$ws1.Range("C15").Value = 98289.71
$ws1.Range("D15").Value = -73.509999999994761
$ws1.Range("E15").Value = -0.00074733218371658694
$ws1.Range("F15").Value = -0.1590395345844621
$ws1.Range("G15").Value = 0.01502988327974286
$ws1.Range("H15").Value = -0.34862410954795853
$ws1.Range("B16").Value = "Indonesia"
$ws1.Range("C16").Value = 5059.22
$ws1.Range("D16").Value = 42.510000000000218
$ws1.Range("E16").Value = 0.0084736809582375283
$ws1.Range("F16").Value = -0.19148208521110999
$ws1.Range("G16").Value = 0.0132132633623927
$ws1.Range("H16").Value = -0.23815728818161699
$ws1.Range("B17").Value = "Japan"
$ws1.Range("C17").Value = 23360.3
$ws1.Range("D17").Value = -46.190000000002328
$ws1.Range("E17").Value = -0.0019733843049514199
$ws1.Range("F17").Value = 0.0066985967594719309
$ws1.Range("G17").Value = 0.011800990786400961
$ws1.Range("H17").Value = 0.038102075037377503
$ws1.Range("B18").Value = "Switzerland"
$ws1.Range("C18").Value = 10539.17
$ws1.Range("D18").Value = 99.649999999999636
$ws1.Range("E18").Value = 0.0095454580287215318
$ws1.Range("F18").Value = -0.01183639447522411
$ws1.Range("G18").Value = 0.01174606697566194
$ws1.Range("H18").Value = 0.056910347399661047
$ws1.Range("B19").Value = "Thailand"
$ws1.Range("C19").Value = 1288.3900000000001
$ws1.Range("D19").Value = 8.4300000000000637
$ws1.Range("E19").Value = 0.0065861433169787098
$ws1.Range("F19").Value = -0.1785846350015938
$ws1.Range("G19").Value = 0.010785692054639419
$ws1.Range("H19").Value = -0.20536762866068781
$ws1.Range("B20").Value = "Singapore"
$ws1.Range("C20").Value = 264.88
$ws1.Range("D20").Value = 0.45999999999997948
$ws1.Range("E20").Value = 0.00173965660691322
$ws1.Range("F20").Value = -0.2114554493763211
$ws1.Range("G20").Value = 0.01067822899868465
$ws1.Range("H20").Value = -0.21573764659891301
$ws1.Range("B21").Value = "United Kingdom"
$ws1.Range("C21").Value = 6007.05
$ws1.Range("D21").Value = -25.03999999999996
$ws1.Range("E21").Value = -0.0041511316973055212
$ws1.Range("F21").Value = -0.20702569125610201
$ws1.Range("G21").Value = 0.008853304685087382
$ws1.Range("H21").Value = -0.21336384288413451
$ws1.Range("B22").Value = "Australia"
$ws1.Range("C22").Value = 5864.5
$ws1.Range("D22").Value = 5.1000000000003638
$ws1.Range("E22").Value = 0.00087039628630924959
$ws1.Range("F22").Value = -0.12934067728669621
$ws1.Range("G22").Value = 0.0077726556149735337
$ws1.Range("H22").Value = -0.083089189431478694
$ws1.Range("B23").Value = "New Zealand"
$ws1.Range("C23").Value = 11633.52
$ws1.Range("D23").Value = -114.5100000000002
$ws1.Range("E23").Value = -0.0097471661206176563
$ws1.Range("F23").Value = 0.00053322691729484895
$ws1.Range("G23").Value = 0.0074902166909374124
$ws1.Range("H23").Value = 0.01667759750142506
$ws1.Range("B24").Value = "Malaysia"
$ws1.Range("C24").Value = 1506.63
$ws1.Range("D24").Value = 1.7800000000002001
$ws1.Range("E24").Value = 0.0011828421437354299
$ws1.Range("F24").Value = -0.057036100540757029
$ws1.Range("G24").Value = 0.0072330238369784361
$ws1.Range("H24").Value = -0.065469708389487069
$ws1.Range("B25").Value = "India"
$ws1.Range("C25").Value = 11504.95
$ws1.Range("D25").Value = 40.5
$ws1.Range("E25").Value = 0.0035326596565905528
$ws1.Range("F25").Value = -0.040698571255852189
$ws1.Range("G25").Value = 0.0027647188103301361
$ws1.Range("H25").Value = -0.064850732598673799
$ws1.Range("B26").Value = "Canada"
$ws1.Range("C26").Value = 16198.97
$ws1.Range("D26").Value = -23.489999999999782
$ws1.Range("E26").Value = -0.001447992474630788
$ws1.Range("F26").Value = -0.052994743786636822
$ws1.Range("G26").Value = 0.0015362501837443701
$ws1.Range("H26").Value = -0.064776317107975601
$ws1.Range("B27").Value = "Spain"
$ws1.Range("C27").Value = 6929.8
$ws1.Range("D27").Value = -13.39999999999964
$ws1.Range("E27").Value = -0.0019299458462955781
$ws1.Range("F27").Value = -0.27821350081763158
$ws1.Range("G27").Value = 0.000080278713535442137
$ws1.Range("H27").Value = -0.23383788106016271
$ws1.Range("B28").Value = "Netherlands"
$ws1.Range("C28").Value = 550.85
$ws1.Range("D28").Value = -1.1599999999999679
$ws1.Range("E28").Value = -0.0021014112063186818
$ws1.Range("F28").Value = -0.092862789012581493
$ws1.Range("G28").Value = -0.000091531996873173149
$ws1.Range("H28").Value = -0.037091759784092522
$ws1.Range("B29").Value = "India"
$ws1.Range("C29").Value = 38845.82
$ws1.Range("D29").Value = -8.7300000000032014
$ws1.Range("E29").Value = -0.00022468411035525551
$ws1.Range("F29").Value = -0.045008890854527417
$ws1.Range("G29").Value = -0.00098974969618015329
$ws1.Range("H29").Value = -0.069052532047928827
$ws1.Range("B30").Value = "Hong Kong"
$ws1.Range("C30").Value = 24455.41
$ws1.Range("D30").Value = -47.900000000001462
$ws1.Range("E30").Value = -0.001954837938221488
$ws1.Range("F30").Value = -0.1335915332533367
$ws1.Range("G30").Value = -0.0019162254300245649
$ws1.Range("H30").Value = -0.130058751927729
$ws1.Range("B31").Value = "Peru"
$ws1.Range("C31").Value = 17907.13
$ws1.Range("D31").Value = -118.8099999999977
$ws1.Range("E31").Value = -0.0065910571099203619
$ws1.Range("F31").Value = -0.1293859169337026
$ws1.Range("G31").Value = -0.003498128508283127
$ws1.Range("H31").Value = -0.1833526906492268
$ws1.Range("B32").Value = "Germany"
$ws1.Range("C32").Value = 13116.25
$ws1.Range("D32").Value = -86.590000000000146
$ws1.Range("E32").Value = -0.0065584374270990509
$ws1.Range("F32").Value = -0.00081816166539316892
$ws1.Range("G32").Value = -0.0045575351662510144
$ws1.Range("H32").Value = 0.06061179494467428
$ws1.Range("B33").Value = "United States"
$ws1.Range("C33").Value = 10793.28
$ws1.Range("D33").Value = -60.260000000000218
$ws1.Range("E33").Value = -0.0055521055802991706
$ws1.Range("F33").Value = 0.18980495994585239
$ws1.Range("G33").Value = -0.0055521055802991706
$ws1.Range("H33").Value = 0.18980495994585239
$ws1.Range("B34").Value = "Colombia"
$ws1.Range("C34").Value = 1206.55
$ws1.Range("D34").Value = -5.9000000000000909
$ws1.Range("E34").Value = -0.0048661800486619064
$ws1.Range("F34").Value = -0.280311842003233
$ws1.Range("G34").Value = -0.0062383109016326266
$ws1.Range("H34").Value = -0.3687802776537018
$ws1.Range("B35").Value = "United States"
$ws1.Range("C35").Value = 3319.47
$ws1.Range("D35").Value = -21.5
$ws1.Range("E35").Value = -0.0064352568266102814
$ws1.Range("F35").Value = 0.022545806276722761
$ws1.Range("G35").Value = -0.0064352568266102814
$ws1.Range("H35").Value = 0.022545806276722761
$ws1.Range("B36").Value = "Philippines"
$ws1.Range("C36").Value = 5908.9
$ws1.Range("D36").Value = -59.0600000000004
$ws1.Range("E36").Value = -0.0098961789288132262
$ws1.Range("F36").Value = -0.24224179166875059
$ws1.Range("G36").Value = -0.0072343354986501573
$ws1.Range("H36").Value = -0.1997427405484794
$ws1.Range("B37").Value = "Germany"
$ws1.Range("C37").Value = 3283.69
$ws1.Range("D37").Value = -32.119999999999891
$ws1.Range("E37").Value = -0.0096869241603106948
$ws1.Range("F37").Value = -0.1249373754170531
$ws1.Range("G37").Value = -0.0076923230211525251
$ws1.Range("H37").Value = -0.071138300016704736
$ws1.Range("B38").Value = "France"
$ws1.Range("C38").Value = 4978.18
$ws1.Range("D38").Value = -55.960000000000043
$ws1.Range("E38").Value = -0.01111609927415602
$ws1.Range("F38").Value = -0.17217834937200571
$ws1.Range("G38").Value = -0.0091243766532963511
$ws1.Range("H38").Value = -0.1212836612102275
$ws1.Range("B39").Value = "Italy"
$ws1.Range("C39").Value = 19524.939999999999
$ws1.Range("D39").Value = -295.81000000000131
$ws1.Range("E39").Value = -0.014924258668314819
$ws1.Range("F39").Value = -0.17201561068117999
$ws1.Range("G39").Value = -0.01294020610579427
$ws1.Range("H39").Value = -0.12111091730280089

# --- Commodities sheet: update row 181 precious metals close + add rows 182-186 ---
$ws2.Range("K181").Value = 192.15
$ws2.Range("L181").Value = 147.25
$ws2.Range("M181").Value = 1111.7

$ws2.Range("A182").Value = 44088
$ws2.Range("B182").Value = 1963.7
$ws2.Range("C182").Value = 27.355
$ws2.Range("D182").Value = 965
$ws2.Range("E182").Value = 2328.1
$ws2.Range("F182").Value = 37.26
$ws2.Range("G182").Value = 39.61
$ws2.Range("H182").Value = 2.31
$ws2.Range("I182").Value = 1.0933999999999999
$ws2.Range("J182").Value = 1778.01
$ws2.Range("K182").Value = 192.25
$ws2.Range("L182").Value = 149.30000000000001
$ws2.Range("M182").Value = 1121.4000000000001
$ws2.Range("N182").Value = 3.0680000000000001
$ws2.Range("O182").Value = 357.5
$ws2.Range("P182").Value = 1002.12
$ws2.Range("Q182").Value = 66.650000000000006
$ws2.Range("R182").Value = 2628
$ws2.Range("S182").Value = 123.55
$ws2.Range("T182").Value = 11.76
$ws2.Range("U182").Value = 106.875
$ws2.Range("V182").Value = 64.63
$ws2.Range("W182").Value = 984.5

$ws2.Range("A183").Value = 44089
$ws2.Range("B183").Value = 1966.2
$ws2.Range("C183").Value = 27.463999999999999
$ws2.Range("D183").Value = 975.7
$ws2.Range("E183").Value = 2426.8000000000002
$ws2.Range("F183").Value = 38.28
$ws2.Range("G183").Value = 40.53
$ws2.Range("H183").Value = 2.3620000000000001
$ws2.Range("I183").Value = 1.0992999999999999
$ws2.Range("J183").Value = 1793.5
$ws2.Range("K183").Value = 194.3
$ws2.Range("L183").Value = 148.69999999999999
$ws2.Range("M183").Value = 1119.5999999999999
$ws2.Range("N183").Value = 3.0630000000000002
$ws2.Range("O183").Value = 366
$ws2.Range("P183").Value = 989.38
$ws2.Range("Q183").Value = 66.37
$ws2.Range("R183").Value = 2766
$ws2.Range("S183").Value = 122.2
$ws2.Range("T183").Value = 12.08
$ws2.Range("U183").Value = 107.1
$ws2.Range("V183").Value = 65.7
$ws2.Range("W183").Value = 920

$ws2.Range("A184").Value = 44090
$ws2.Range("B184").Value = 1970.5
$ws2.Range("C184").Value = 27.475999999999999
$ws2.Range("D184").Value = 972.05
$ws2.Range("E184").Value = 2413.9
$ws2.Range("F184").Value = 40.159999999999997
$ws2.Range("G184").Value = 42.22
$ws2.Range("H184").Value = 2.2669999999999999
$ws2.Range("I184").Value = 1.1163000000000001
$ws2.Range("J184").Value = 1772.35
$ws2.Range("K184").Value = 195.35
$ws2.Range("L184").Value = 147.5
$ws2.Range("M184").Value = 1115.8
$ws2.Range("N184").Value = 3.0609999999999999
$ws2.Range("O184").Value = 371.75
$ws2.Range("P184").Value = 1011.62
$ws2.Range("Q184").Value = 66.489999999999995
$ws2.Range("R184").Value = 2602
$ws2.Range("S184").Value = 120.55
$ws2.Range("T184").Value = 12.35
$ws2.Range("U184").Value = 106.72499999999999
$ws2.Range("V184").Value = 65.22
$ws2.Range("W184").Value = 612

$ws2.Range("A185").Value = 44091
$ws2.Range("B185").Value = 1949.9
$ws2.Range("C185").Value = 27.1
$ws2.Range("D185").Value = 937.6
$ws2.Range("E185").Value = 2349.1
$ws2.Range("F185").Value = 40.97
$ws2.Range("G185").Value = 43.3
$ws2.Range("H185").Value = 2.0419999999999998
$ws2.Range("I185").Value = 1.1597999999999999
$ws2.Range("J185").Value = 1757
$ws2.Range("K185").Value = 195.3
$ws2.Range("L185").Value = 148.1
$ws2.Range("M185").Value = 1107.7
$ws2.Range("N185").Value = 3.0710000000000002
$ws2.Range("O185").Value = 375.25
$ws2.Range("P185").Value = 1027.75
$ws2.Range("Q185").Value = 65.84
$ws2.Range("R185").Value = 2600
$ws2.Range("S185").Value = 118.1
$ws2.Range("T185").Value = 12.62
$ws2.Range("U185").Value = 106.77500000000001
$ws2.Range("V185").Value = 66.53
$ws2.Range("W185").Value = 597.6

$ws2.Range("A186").Value = 44092
$ws2.Range("B186").Value = 1962.1
$ws2.Range("C186").Value = 27.129000000000001
$ws2.Range("D186").Value = 932.4
$ws2.Range("E186").Value = 2382.6
$ws2.Range("F186").Value = 41.11
$ws2.Range("G186").Value = 43.15
$ws2.Range("H186").Value = 2.048
$ws2.Range("I186").Value = 1.159
$ws2.Range("J186").Value = 1795.5
$ws2.Range("K186").Value = 197.85
$ws2.Range("L186").Value = 148.4
$ws2.Range("M186").Value = 1090.8
$ws2.Range("N186").Value = 3.1160000000000001
$ws2.Range("O186").Value = 378.5
$ws2.Range("P186").Value = 1043.3800000000001
$ws2.Range("Q186").Value = 65.62
$ws2.Range("R186").Value = 2641
$ws2.Range("S186").Value = 113.6
$ws2.Range("T186").Value = 12.77
$ws2.Range("U186").Value = 107.35
$ws2.Range("V186").Value = 66.5
$ws2.Range("W186").Value = 578.6

# --- view state: activate Commodities sheet, select all cells, set workbook active tab ---
$ws2.Activate()
$ws2.Cells.Select()